# Natmi following Dr Hou advice
# Update LR-pair edge-weight stats for Adam17-Itgb1: the number of ligand-
# and receptor-expressing cells moves from 1 to 3 for every row, which
# changes the dependent total/specificity metrics across columns G-T.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 30.07831966666667
$ws.Range("H2").Value = 90.234959
$ws.Range("I2").Value = 0.2269842729019557
$ws.Range("J2").Value = 0.2269842729019557
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 98.946724
$ws.Range("N2").Value = 296.840172
$ws.Range("O2").Value = 0.2098009692989996
$ws.Range("P2").Value = 0.2098009692989996
$ws.Range("Q2").Value = 2976.151194441439
$ws.Range("R2").Value = 26785.36074997295
$ws.Range("S2").Value = 0.04762152047045896
$ws.Range("T2").Value = 0.04762152047045896

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 30.07831966666667
$ws.Range("H3").Value = 90.234959
$ws.Range("I3").Value = 0.2269842729019557
$ws.Range("J3").Value = 0.2269842729019557
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 163.0062356666667
$ws.Range("N3").Value = 489.018707
$ws.Range("O3").Value = 0.345629090707923
$ws.Range("P3").Value = 0.3456290907079231
$ws.Range("Q3").Value = 4902.953664042002
$ws.Range("R3").Value = 44126.58297637801
$ws.Range("S3").Value = 0.07845236784810201
$ws.Range("T3").Value = 0.07845236784810201

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 30.07831966666667
$ws.Range("H4").Value = 90.234959
$ws.Range("I4").Value = 0.2269842729019557
$ws.Range("J4").Value = 0.2269842729019557
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 65.39610666666668
$ws.Range("N4").Value = 196.18832
$ws.Range("O4").Value = 0.1386621609326595
$ws.Range("P4").Value = 0.1386621609326595
$ws.Range("Q4").Value = 1967.005001275432
$ws.Range("R4").Value = 17703.04501147888
$ws.Range("S4").Value = 0.03147412977831367
$ws.Range("T4").Value = 0.03147412977831368

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 30.07831966666667
$ws.Range("H5").Value = 90.234959
$ws.Range("I5").Value = 0.2269842729019557
$ws.Range("J5").Value = 0.2269842729019557
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 144.2727966666667
$ws.Range("N5").Value = 432.81839
$ws.Range("O5").Value = 0.3059077790604178
$ws.Range("P5").Value = 0.3059077790604179
$ws.Range("Q5").Value = 4339.483297344002
$ws.Range("R5").Value = 39055.34967609601
$ws.Range("S5").Value = 0.06943625480508106
$ws.Range("T5").Value = 0.06943625480508107

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 31.96959266666667
$ws.Range("H6").Value = 95.90877800000001
$ws.Range("I6").Value = 0.2412566535243296
$ws.Range("J6").Value = 0.2412566535243296
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 98.946724
$ws.Range("N6").Value = 296.840172
$ws.Range("O6").Value = 0.2098009692989996
$ws.Range("P6").Value = 0.2098009692989996
$ws.Range("Q6").Value = 3163.286461981091
$ws.Range("R6").Value = 28469.57815782982
$ws.Range("S6").Value = 0.05061587975923726
$ws.Range("T6").Value = 0.05061587975923727

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 31.96959266666667
$ws.Range("H7").Value = 95.90877800000001
$ws.Range("I7").Value = 0.2412566535243296
$ws.Range("J7").Value = 0.2412566535243296
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 163.0062356666667
$ws.Range("N7").Value = 489.018707
$ws.Range("O7").Value = 0.345629090707923
$ws.Range("P7").Value = 0.3456290907079231
$ws.Range("Q7").Value = 5211.242956390006
$ws.Range("R7").Value = 46901.18660751006
$ws.Range("S7").Value = 0.08338531778485048
$ws.Range("T7").Value = 0.08338531778485049

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 31.96959266666667
$ws.Range("H8").Value = 95.90877800000001
$ws.Range("I8").Value = 0.2412566535243296
$ws.Range("J8").Value = 0.2412566535243296
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 65.39610666666668
$ws.Range("N8").Value = 196.18832
$ws.Range("O8").Value = 0.1386621609326595
$ws.Range("P8").Value = 0.1386621609326595
$ws.Range("Q8").Value = 2090.686892119219
$ws.Range("R8").Value = 18816.18202907297
$ws.Range("S8").Value = 0.03345316891706546
$ws.Range("T8").Value = 0.03345316891706546

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 31.96959266666667
$ws.Range("H9").Value = 95.90877800000001
$ws.Range("I9").Value = 0.2412566535243296
$ws.Range("J9").Value = 0.2412566535243296
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 144.2727966666667
$ws.Range("N9").Value = 432.81839
$ws.Range("O9").Value = 0.3059077790604178
$ws.Range("P9").Value = 0.3059077790604179
$ws.Range("Q9").Value = 4612.342542314158
$ws.Range("R9").Value = 41511.08288082743
$ws.Range("S9").Value = 0.0738022870631764
$ws.Range("T9").Value = 0.07380228706317642

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 55.83720666666667
$ws.Range("H10").Value = 167.51162
$ws.Range("I10").Value = 0.4213722008598541
$ws.Range("J10").Value = 0.4213722008598541
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 98.946724
$ws.Range("N10").Value = 296.840172
$ws.Range("O10").Value = 0.2098009692989996
$ws.Range("P10").Value = 0.2098009692989996
$ws.Range("Q10").Value = 5524.908676977627
$ws.Range("R10").Value = 49724.17809279864
$ws.Range("S10").Value = 0.08840429617605014
$ws.Range("T10").Value = 0.08840429617605015

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 55.83720666666667
$ws.Range("H11").Value = 167.51162
$ws.Range("I11").Value = 0.4213722008598541
$ws.Range("J11").Value = 0.4213722008598541
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 163.0062356666667
$ws.Range("N11").Value = 489.018707
$ws.Range("O11").Value = 0.345629090707923
$ws.Range("P11").Value = 0.3456290907079231
$ws.Range("Q11").Value = 9101.812868875038
$ws.Range("R11").Value = 81916.31581987534
$ws.Range("S11").Value = 0.1456384906327877
$ws.Range("T11").Value = 0.1456384906327877

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 55.83720666666667
$ws.Range("H12").Value = 167.51162
$ws.Range("I12").Value = 0.4213722008598541
$ws.Range("J12").Value = 0.4213722008598541
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 65.39610666666668
$ws.Range("N12").Value = 196.18832
$ws.Range("O12").Value = 0.1386621609326595
$ws.Range("P12").Value = 0.1386621609326595
$ws.Range("Q12").Value = 3651.535923142045
$ws.Range("R12").Value = 32863.82330827841
$ws.Range("S12").Value = 0.058428379928178
$ws.Range("T12").Value = 0.05842837992817801

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 55.83720666666667
$ws.Range("H13").Value = 167.51162
$ws.Range("I13").Value = 0.4213722008598541
$ws.Range("J13").Value = 0.4213722008598541
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 144.2727966666667
$ws.Range("N13").Value = 432.81839
$ws.Range("O13").Value = 0.3059077790604178
$ws.Range("P13").Value = 0.3059077790604179
$ws.Range("Q13").Value = 8055.789963854644
$ws.Range("R13").Value = 72502.1096746918
$ws.Range("S13").Value = 0.1289010341228383
$ws.Range("T13").Value = 0.1289010341228383

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 14.62767266666667
$ws.Range("H14").Value = 43.883018
$ws.Range("I14").Value = 0.1103868727138606
$ws.Range("J14").Value = 0.1103868727138606
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 98.946724
$ws.Range("N14").Value = 296.840172
$ws.Range("O14").Value = 0.2098009692989996
$ws.Range("P14").Value = 0.2098009692989996
$ws.Range("Q14").Value = 1447.360290111011
$ws.Range("R14").Value = 13026.2426109991
$ws.Range("S14").Value = 0.02315927289325325
$ws.Range("T14").Value = 0.02315927289325325

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 14.62767266666667
$ws.Range("H15").Value = 43.883018
$ws.Range("I15").Value = 0.1103868727138606
$ws.Range("J15").Value = 0.1103868727138606
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 163.0062356666667
$ws.Range("N15").Value = 489.018707
$ws.Range("O15").Value = 0.345629090707923
$ws.Range("P15").Value = 0.3456290907079231
$ws.Range("Q15").Value = 2384.401857957525
$ws.Range("R15").Value = 21459.61672161773
$ws.Range("S15").Value = 0.0381529144421829
$ws.Range("T15").Value = 0.0381529144421829

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 14.62767266666667
$ws.Range("H16").Value = 43.883018
$ws.Range("I16").Value = 0.1103868727138606
$ws.Range("J16").Value = 0.1103868727138606
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 65.39610666666668
$ws.Range("N16").Value = 196.18832
$ws.Range("O16").Value = 0.1386621609326595
$ws.Range("P16").Value = 0.1386621609326595
$ws.Range("Q16").Value = 956.592841994418
$ws.Range("R16").Value = 8609.335577949762
$ws.Range("S16").Value = 0.01530648230910234
$ws.Range("T16").Value = 0.01530648230910234

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 14.62767266666667
$ws.Range("H17").Value = 43.883018
$ws.Range("I17").Value = 0.1103868727138606
$ws.Range("J17").Value = 0.1103868727138606
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 144.2727966666667
$ws.Range("N17").Value = 432.81839
$ws.Range("O17").Value = 0.3059077790604178
$ws.Range("P17").Value = 0.3059077790604179
$ws.Range("Q17").Value = 2110.375244344558
$ws.Range("R17").Value = 18993.37719910102
$ws.Range("S17").Value = 0.03376820306932215
$ws.Range("T17").Value = 0.03376820306932216
